# Update the cryptos list (Price and Volume(1h) columns) with refreshed
# values, as performed by the scheduled GitHub Actions job.
#
# Note: some "Price" values (column D) look like plain decimals (e.g.
# "23.00"), which Excel's COM Value setter would otherwise auto-convert
# to a numeric type and silently drop the trailing zero. To keep those
# cells as text (matching the workbook's original inlineStr formatting)
# we prefix them with an apostrophe to force text entry, then reset the
# cell style back to "Normal" so no quotePrefix style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.213.32"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.275.49"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'588.04"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'185.90"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.134"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "3.848.10"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "'28.80"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "68.219.10"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "3.268.60"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "'13.68"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "'383.43"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'0.517"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("D26").Value = "'9.94"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "'0.184"
$ws.Range("E27").Value = "  +3.38%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'5.82"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'7.27"
$ws.Range("E31").Value = "  +6.53%  "
$ws.Range("D32").Value = "'23.00"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").Value = "'162.81"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "'26.86"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("E41").Value = "  +6.01%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").Value = "'350.21"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'41.55"
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").Value = "2.658.66"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").Value = "'32.22"
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +0.16%  "

# Reset style on cells where we used an apostrophe-prefix to force text,
# so no quotePrefix style attribute lingers on the cell (matches original formatting).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D49").Style = "Normal"
